$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format temporarily so numeric-looking price strings
# (e.g. "558.07") are preserved as text, matching the source data which stores
# all values as strings. Style is reset back to Normal afterwards so no cell
# styling changes remain in the output.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '60.374.42'
$ws.Range("E2").Value = '  -3.09%  '
$ws.Range("D3").Value = '3.305.10'
$ws.Range("E3").Value = '  -3.66%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '558.07'
$ws.Range("E5").Value = '  -3.67%  '
$ws.Range("D6").Value = '142.22'
$ws.Range("E6").Value = '  -6.91%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.304.97'
$ws.Range("E8").Value = '  -3.66%  '
$ws.Range("D9").Value = '0.468'
$ws.Range("E9").Value = '  -3.02%  '
$ws.Range("E10").Value = '  -2.40%  '
$ws.Range("E11").Value = '  -4.75%  '
$ws.Range("E12").Value = '  -2.43%  '
$ws.Range("D13").Value = '3.872.78'
$ws.Range("E13").Value = '  -3.48%  '
$ws.Range("E14").Value = '  +0.65%  '
$ws.Range("D15").Value = '26.91'
$ws.Range("E15").Value = '  -6.32%  '
$ws.Range("D16").Value = '3.306.38'
$ws.Range("E16").Value = '  -4.65%  '
$ws.Range("D17").Value = '0.0000165'
$ws.Range("E17").Value = '  -3.90%  '
$ws.Range("D18").Value = '60.390.39'
$ws.Range("E18").Value = '  -3.08%  '
$ws.Range("D19").Value = '6.10'
$ws.Range("E19").Value = '  -6.34%  '
$ws.Range("D20").Value = '13.99'
$ws.Range("E20").Value = '  -3.71%  '
$ws.Range("D21").Value = '8.65'
$ws.Range("E21").Value = '  -3.52%  '
$ws.Range("D22").Value = '375.19'
$ws.Range("E22").Value = '  -2.24%  '
$ws.Range("D23").Value = '74.38'
$ws.Range("E23").Value = '  -1.11%  '
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").Value = '0.535'
$ws.Range("E25").Value = '  -6.33%  '
$ws.Range("D26").Value = '3.437.20'
$ws.Range("E26").Value = '  -3.45%  '
$ws.Range("D27").Value = '0.0000102'
$ws.Range("E27").Value = '  -9.18%  '
$ws.Range("E28").Value = '  -4.34%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").Value = '7.18'
$ws.Range("E30").Value = '  -6.77%  '
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("E32").Value = '  -4.11%  '
$ws.Range("D33").Value = '7.56'
$ws.Range("E33").Value = '  -5.20%  '
$ws.Range("D34").Value = '22.63'
$ws.Range("E34").Value = '  -2.70%  '
$ws.Range("D35").Value = '1.24'
$ws.Range("E35").Value = '  -7.71%  '
$ws.Range("E36").Value = '  -6.55%  '
$ws.Range("D37").Value = '166.90'
$ws.Range("E37").Value = '  -1.19%  '
$ws.Range("D38").Value = '1.53'
$ws.Range("E38").Value = '  -5.53%  '
$ws.Range("D39").Value = '6.69'
$ws.Range("E39").Value = '  -3.75%  '
$ws.Range("D40").Value = '3.338.11'
$ws.Range("E40").Value = '  -3.60%  '
$ws.Range("D41").Value = '26.64'
$ws.Range("E41").Value = '  -14.35%  '
$ws.Range("D42").Value = '0.0730'
$ws.Range("E42").Value = '  -6.91%  '
$ws.Range("D43").Value = '41.93'
$ws.Range("E43").Value = '  -2.09%  '
$ws.Range("E44").Value = '  -3.66%  '
$ws.Range("D45").Value = '4.13'
$ws.Range("E45").Value = '  -6.46%  '
$ws.Range("D46").Value = '1.11'
$ws.Range("E46").Value = '  -5.95%  '
$ws.Range("D47").Value = '1.57'
$ws.Range("E47").Value = '  -6.78%  '
$ws.Range("D48").Value = '2.351.63'
$ws.Range("E48").Value = '  -7.46%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").Value = '6.39'
$ws.Range("E50").Value = '  -7.42%  '
$ws.Range("D51").Value = '21.34'
$ws.Range("E51").Value = '  -5.72%  '

$dRange.Style = "Normal"
